$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 0.05441666666666667
$ws.Cells.Item(2, 8).Value = 0.16325
$ws.Cells.Item(2, 9).Value = 0.00608027172874025
$ws.Cells.Item(2, 10).Value = 0.006080271728740251
$ws.Cells.Item(2, 13).Value = 2.112132333333333
$ws.Cells.Item(2, 14).Value = 6.336397
$ws.Cells.Item(2, 15).Value = 0.09680078109791683
$ws.Cells.Item(2, 16).Value = 0.1001796274503123
$ws.Cells.Item(2, 17).Value = 0.1149352011388889
$ws.Cells.Item(2, 18).Value = 1.03441681025
$ws.Cells.Item(2, 19).Value = 0.0005885750526296373
$ws.Cells.Item(2, 20).Value = 0.0006091193565818649
$ws.Cells.Item(3, 7).Value = 0.05441666666666667
$ws.Cells.Item(3, 8).Value = 0.16325
$ws.Cells.Item(3, 9).Value = 0.00608027172874025
$ws.Cells.Item(3, 10).Value = 0.006080271728740251
$ws.Cells.Item(3, 13).Value = 8.637706333333334
$ws.Cells.Item(3, 15).Value = 0.3958732636044222
$ws.Cells.Item(3, 16).Value = 0.4096912815746252
$ws.Cells.Item(3, 17).Value = 0.4700351863055556
$ws.Cells.Item(3, 18).Value = 4.23031667675
$ws.Cells.Item(3, 19).Value = 0.002407017012858105
$ws.Cells.Item(3, 20).Value = 0.002491034316869555
$ws.Cells.Item(4, 7).Value = 0.05441666666666667
$ws.Cells.Item(4, 8).Value = 0.16325
$ws.Cells.Item(4, 9).Value = 0.00608027172874025
$ws.Cells.Item(4, 10).Value = 0.006080271728740251
$ws.Cells.Item(4, 13).Value = 4.811824666666667
$ws.Cells.Item(4, 14).Value = 14.435474
$ws.Cells.Item(4, 15).Value = 0.2205299255584317
$ws.Cells.Item(4, 16).Value = 0.2282275569836723
$ws.Cells.Item(4, 17).Value = 0.2618434589444445
$ws.Cells.Item(4, 18).Value = 2.3565911305
$ws.Cells.Item(4, 19).Value = 0.001340881871714124
$ws.Cells.Item(4, 20).Value = 0.001387685562447277
$ws.Cells.Item(5, 7).Value = 0.05441666666666667
$ws.Cells.Item(5, 8).Value = 0.16325
$ws.Cells.Item(5, 9).Value = 0.00608027172874025
$ws.Cells.Item(5, 10).Value = 0.006080271728740251
$ws.Cells.Item(5, 13).Value = 2.2077635
$ws.Cells.Item(5, 14).Value = 4.415527
$ws.Cells.Item(5, 15).Value = 0.1011836369846164
$ws.Cells.Item(5, 16).Value = 0.06981031173659025
$ws.Cells.Item(5, 17).Value = 0.1201391304583333
$ws.Cells.Item(5, 18).Value = 0.72083478275
$ws.Cells.Item(5, 19).Value = 0.0006152240073686792
$ws.Cells.Item(5, 20).Value = 0.0004244656648265334
$ws.Cells.Item(6, 7).Value = 0.05441666666666667
$ws.Cells.Item(6, 8).Value = 0.16325
$ws.Cells.Item(6, 9).Value = 0.00608027172874025
$ws.Cells.Item(6, 10).Value = 0.006080271728740251
$ws.Cells.Item(6, 13).Value = 4.049945999999999
$ws.Cells.Item(6, 14).Value = 12.149838
$ws.Cells.Item(6, 15).Value = 0.185612392754613
$ws.Cells.Item(6, 16).Value = 0.1920912222548
$ws.Cells.Item(6, 17).Value = 0.2203845615
$ws.Cells.Item(6, 18).Value = 1.9834610535
$ws.Cells.Item(6, 19).Value = 0.001128573784169705
$ws.Cells.Item(6, 20).Value = 0.001167966828015021
$ws.Cells.Item(7, 7).Value = 4.046611333333334
$ws.Cells.Item(7, 8).Value = 12.139834
$ws.Cells.Item(7, 9).Value = 0.4521500120171497
$ws.Cells.Item(7, 10).Value = 0.4521500120171497
$ws.Cells.Item(7, 13).Value = 2.112132333333333
$ws.Cells.Item(7, 14).Value = 6.336397
$ws.Cells.Item(7, 15).Value = 0.09680078109791683
$ws.Cells.Item(7, 16).Value = 0.1001796274503123
$ws.Cells.Item(7, 17).Value = 8.546978637566445
$ws.Cells.Item(7, 18).Value = 76.922807738098
$ws.Cells.Item(7, 19).Value = 0.04376847433669257
$ws.Cells.Item(7, 20).Value = 0.0452962197555323
$ws.Cells.Item(8, 7).Value = 4.046611333333334
$ws.Cells.Item(8, 8).Value = 12.139834
$ws.Cells.Item(8, 9).Value = 0.4521500120171497
$ws.Cells.Item(8, 10).Value = 0.4521500120171497
$ws.Cells.Item(8, 13).Value = 8.637706333333334
$ws.Cells.Item(8, 15).Value = 0.3958732636044222
$ws.Cells.Item(8, 16).Value = 0.4096912815746252
$ws.Cells.Item(8, 17).Value = 34.95344034247179
$ws.Cells.Item(8, 18).Value = 314.5809630822461
$ws.Cells.Item(8, 19).Value = 0.1789941008960077
$ws.Cells.Item(8, 20).Value = 0.1852419178872882
$ws.Cells.Item(9, 7).Value = 4.046611333333334
$ws.Cells.Item(9, 8).Value = 12.139834
$ws.Cells.Item(9, 9).Value = 0.4521500120171497
$ws.Cells.Item(9, 10).Value = 0.4521500120171497
$ws.Cells.Item(9, 13).Value = 4.811824666666667
$ws.Cells.Item(9, 14).Value = 14.435474
$ws.Cells.Item(9, 15).Value = 0.2205299255584317
$ws.Cells.Item(9, 16).Value = 0.2282275569836723
$ws.Cells.Item(9, 17).Value = 19.47158423014623
$ws.Cells.Item(9, 18).Value = 175.244258071316
$ws.Cells.Item(9, 19).Value = 0.09971260849138602
$ws.Cells.Item(9, 20).Value = 0.1031930926328121
$ws.Cells.Item(10, 7).Value = 4.046611333333334
$ws.Cells.Item(10, 8).Value = 12.139834
$ws.Cells.Item(10, 9).Value = 0.4521500120171497
$ws.Cells.Item(10, 10).Value = 0.4521500120171497
$ws.Cells.Item(10, 13).Value = 2.2077635
$ws.Cells.Item(10, 14).Value = 4.415527
$ws.Cells.Item(10, 15).Value = 0.1011836369846164
$ws.Cells.Item(10, 16).Value = 0.06981031173659025
$ws.Cells.Item(10, 17).Value = 8.933960800419667
$ws.Cells.Item(10, 18).Value = 53.603764802518
$ws.Cells.Item(10, 19).Value = 0.0457501826785332
$ws.Cells.Item(10, 20).Value = 0.03156473329062024
$ws.Cells.Item(11, 7).Value = 4.046611333333334
$ws.Cells.Item(11, 8).Value = 12.139834
$ws.Cells.Item(11, 9).Value = 0.4521500120171497
$ws.Cells.Item(11, 10).Value = 0.4521500120171497
$ws.Cells.Item(11, 13).Value = 4.049945999999999
$ws.Cells.Item(11, 14).Value = 12.149838
$ws.Cells.Item(11, 15).Value = 0.185612392754613
$ws.Cells.Item(11, 16).Value = 0.1920912222548
$ws.Cells.Item(11, 17).Value = 16.388557382988
$ws.Cells.Item(11, 18).Value = 147.497016446892
$ws.Cells.Item(11, 19).Value = 0.08392464561453017
$ws.Cells.Item(11, 20).Value = 0.08685404845089678
$ws.Cells.Item(12, 7).Value = 4.848681666666667
$ws.Cells.Item(12, 8).Value = 14.546045
$ws.Cells.Item(12, 9).Value = 0.5417697162541101
$ws.Cells.Item(12, 10).Value = 0.5417697162541102
$ws.Cells.Item(12, 13).Value = 2.112132333333333
$ws.Cells.Item(12, 14).Value = 6.336397
$ws.Cells.Item(12, 15).Value = 0.09680078109791683
$ws.Cells.Item(12, 16).Value = 0.1001796274503123
$ws.Cells.Item(12, 17).Value = 10.24105732220722
$ws.Cells.Item(12, 18).Value = 92.16951589986499
$ws.Cells.Item(12, 19).Value = 0.05244373170859462
$ws.Cells.Item(12, 20).Value = 0.05427428833819818
$ws.Cells.Item(13, 7).Value = 4.848681666666667
$ws.Cells.Item(13, 8).Value = 14.546045
$ws.Cells.Item(13, 9).Value = 0.5417697162541101
$ws.Cells.Item(13, 10).Value = 0.5417697162541102
$ws.Cells.Item(13, 13).Value = 8.637706333333334
$ws.Cells.Item(13, 15).Value = 0.3958732636044222
$ws.Cells.Item(13, 16).Value = 0.4096912815746252
$ws.Cells.Item(13, 17).Value = 41.88148834048389
$ws.Cells.Item(13, 18).Value = 376.933395064355
$ws.Cells.Item(13, 19).Value = 0.2144721456955563
$ws.Cells.Item(13, 20).Value = 0.2219583293704674
$ws.Cells.Item(14, 7).Value = 4.848681666666667
$ws.Cells.Item(14, 8).Value = 14.546045
$ws.Cells.Item(14, 9).Value = 0.5417697162541101
$ws.Cells.Item(14, 10).Value = 0.5417697162541102
$ws.Cells.Item(14, 13).Value = 4.811824666666667
$ws.Cells.Item(14, 14).Value = 14.435474
$ws.Cells.Item(14, 15).Value = 0.2205299255584317
$ws.Cells.Item(14, 16).Value = 0.2282275569836723
$ws.Cells.Item(14, 17).Value = 23.33100604448111
$ws.Cells.Item(14, 18).Value = 209.97905440033
$ws.Cells.Item(14, 19).Value = 0.1194764351953316
$ws.Cells.Item(14, 20).Value = 0.1236467787884129
$ws.Cells.Item(15, 7).Value = 4.848681666666667
$ws.Cells.Item(15, 8).Value = 14.546045
$ws.Cells.Item(15, 9).Value = 0.5417697162541101
$ws.Cells.Item(15, 10).Value = 0.5417697162541102
$ws.Cells.Item(15, 13).Value = 2.2077635
$ws.Cells.Item(15, 14).Value = 4.415527
$ws.Cells.Item(15, 15).Value = 0.1011836369846164
$ws.Cells.Item(15, 16).Value = 0.06981031173659025
$ws.Cells.Item(15, 17).Value = 10.70474240678583
$ws.Cells.Item(15, 18).Value = 64.228454440715
$ws.Cells.Item(15, 19).Value = 0.05481823029871448
$ws.Cells.Item(15, 20).Value = 0.03782111278114347
$ws.Cells.Item(16, 7).Value = 4.848681666666667
$ws.Cells.Item(16, 8).Value = 14.546045
$ws.Cells.Item(16, 9).Value = 0.5417697162541101
$ws.Cells.Item(16, 10).Value = 0.5417697162541102
$ws.Cells.Item(16, 13).Value = 4.049945999999999
$ws.Cells.Item(16, 14).Value = 12.149838
$ws.Cells.Item(16, 15).Value = 0.185612392754613
$ws.Cells.Item(16, 16).Value = 0.1920912222548
$ws.Cells.Item(16, 17).Value = 19.63689892119
$ws.Cells.Item(16, 18).Value = 176.73209029071
$ws.Cells.Item(16, 19).Value = 0.1005591733559131
$ws.Cells.Item(16, 20).Value = 0.1040692069758882
